$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "value" header to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Extend the date-formatted style (currently only on A2) down through A22
# before filling in values, so the new rows pick up the same number
# format / font / border as the existing date column.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Date (column A, Excel serial numbers) and year-over-year value
# (column B) series. Column B is blank on the first and last rows
# since there is no prior/next year to compare against.
$dates = @(38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657,46022)
$values = @($null,1.176843378132464,1.383039815128395,-0.6300631236164866,-2.845830838597474,1.163890860292871,1.237492433423526,1.175130261101254,0.09561723522806265,-0.4807826571170737,-0.2020167505668247,-0.3527529999609147,0.2478074346218495,-0.4273761665070541,-0.6791462188813879,-2.118861353231827,1.553047647471506,-0.400533798485958,0.6918534271163068,0.2112326235108375,$null)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    if ($null -eq $values[$i]) {
        $ws.Cells.Item($row, 2).ClearContents() | Out-Null
    } else {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}
